# "fixed allowed stress values and adjusted a, b, t accordingly"
#
# Sheet1 layout:
#   Q18 = a  (label P18="a")
#   Q19 = b  (label P19="b")
#   Q20 = t  (label P20="t")
#   U25 = sigma_dopV (allowed bending stress), label T25="sigma_dopV"
#   U26 = sigma_dopT (allowed torsion/shear stress), label T26="sigma_dopT"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

# Allowed stress values corrected
$ws.Range("U25").Value = 60
$ws.Range("U26").Value = 160

# a, b, t adjusted accordingly
$ws.Range("Q18").Value = 110
$ws.Range("Q19").Value = 90
$ws.Range("Q20").Value = 25
